$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.957.20"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.739.72"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'311.33"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4993"
$ws.Range("E7").Value = "  +8.31%  "
$ws.Range("D8").Value = "'0.3564"
$ws.Range("E8").Value = "  +3.49%  "
$ws.Range("D9").Value = "'42.16"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").Value = "'0.07247"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "'1.060"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'20.29"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "'5.947"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "1.741.60"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'6.836"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "'86.73"
$ws.Range("E17").Value = "  -3.12%  "
$ws.Range("D18").Value = "'0.00001035"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'0.06394"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'16.53"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'5.713"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").Value = "27.018.99"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "'11.38"
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("D25").Value = "'2.044"
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("D26").Value = "'154.80"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'19.80"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "1.941.65"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'2.204"
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").Value = "'120.09"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "'1.045"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").Value = "'0.09480"
$ws.Range("E32").Value = "  +4.21%  "
$ws.Range("D33").Value = "'3.584"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'5.362"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "'0.02197"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "'0.05861"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "'1.429"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").Value = "'0.1997"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'4.769"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").Value = "'0.6023"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").Value = "'1.100"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").Value = "'7.636"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("D44").Value = "'12.78"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "'3.600"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'0.5658"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "'120.21"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "'1.852"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.103"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06671"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  +0.04%  "
